# Updates historico_taxas.xlsx: every existing "Vencimento" row is paired
# with a twin row (same date, the other Taxa flag, a 13:06:11 save time),
# expanding the 14-row history into 28 rows (A1:C29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45792, 0, "2025-04-04 13:05:55"),
    @(45792, 1, "2025-04-04 13:06:11"),
    @(46249, 0, "2025-04-04 13:05:55"),
    @(46249, 1, "2025-04-04 13:06:11"),
    @(46522, 0, "2025-04-04 13:05:55"),
    @(46522, 1, "2025-04-04 13:06:11"),
    @(46980, 0, "2025-04-04 13:05:55"),
    @(46980, 1, "2025-04-04 13:06:11"),
    @(47253, 0, "2025-04-04 13:05:55"),
    @(47253, 1, "2025-04-04 13:06:11"),
    @(47710, 0, "2025-04-04 13:05:55"),
    @(47710, 1, "2025-04-04 13:06:11"),
    @(48441, 0, "2025-04-04 13:05:55"),
    @(48441, 1, "2025-04-04 13:06:11"),
    @(48714, 1, "2025-04-04 13:06:11"),
    @(48714, 0, "2025-04-04 13:05:55"),
    @(49444, 0, "2025-04-04 13:05:55"),
    @(49444, 1, "2025-04-04 13:06:11"),
    @(51363, 1, "2025-04-04 13:06:11"),
    @(51363, 0, "2025-04-04 13:05:55"),
    @(53097, 1, "2025-04-04 13:06:11"),
    @(53097, 0, "2025-04-04 13:05:55"),
    @(55015, 1, "2025-04-04 13:06:11"),
    @(55015, 0, "2025-04-04 13:05:55"),
    @(56749, 1, "2025-04-04 13:06:11"),
    @(56749, 0, "2025-04-04 13:05:55"),
    @(58668, 0, "2025-04-04 13:05:55"),
    @(58668, 1, "2025-04-04 13:06:11"),
)

$rowIndex = 2
foreach ($entry in $data) {
    $dateSerial = $entry[0]
    $taxa = $entry[1]
    $savedAt = $entry[2]

    $cellA = $ws.Cells.Item($rowIndex, 1)
    $cellA.Value = $dateSerial
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($rowIndex, 2).Value = $taxa
    $ws.Cells.Item($rowIndex, 3).Value = $savedAt

    $rowIndex++
}

